# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Row 55: Melilla - Recuperados 60 -> 61, Muertes 2 -> 1
$ws.Range("D55").Value = 61
$ws.Range("E55").Value = 1

# Row 56: Igualada, Vilanova del Cami, Santa Margarida de Montbui y Odena
# Casos totales 58 -> 62, Recuperados 58 -> 60, Muertes 3 -> 2
$ws.Range("B56").Value = 62
$ws.Range("D56").Value = 60
$ws.Range("E56").Value = 2

# Row 57: Ceuta - Casos totales 51 -> 58, Recuperados 50 -> 58, Muertes 1 -> 3
$ws.Range("B57").Value = 58
$ws.Range("D57").Value = 58
$ws.Range("E57").Value = 3

# Update "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 20:20"
